$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ethnicities")

# Clear the previous "Race" table (A1:F7) so stale cells/shared strings are dropped
$ws.Range("A1:G7").Clear()

# New header row (race/ethnicity categories), now starting at column B
$ws.Range("B1").Value = "White"
$ws.Range("C1").Value = "African America"
$ws.Range("D1").Value = "Latinx"
$ws.Range("E1").Value = "Asian"
$ws.Range("F1").Value = "Other"
$ws.Range("G1").Value = "Unknown"

# New row labels (metrics), now in column A
$ws.Range("A2").Value = "Cases"
$ws.Range("A3").Value = "Hospitalizations"
$ws.Range("A4").Value = "ICU"
$ws.Range("A5").Value = "Deaths"

# Cases row
$ws.Range("B2").Value = 1329
$ws.Range("C2").Value = 926
$ws.Range("D2").Value = 291
$ws.Range("E2").Value = 58
$ws.Range("F2").Value = 117

# Hospitalizations row
$ws.Range("B3").Value = 267
$ws.Range("C3").Value = 193
$ws.Range("D3").Value = 53
$ws.Range("E3").Value = 6
$ws.Range("F3").Value = 3

# ICU row
$ws.Range("B4").Value = 95
$ws.Range("C4").Value = 79
$ws.Range("D4").Value = 19
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = 1

# Deaths row
$ws.Range("B5").Value = 150
$ws.Range("C5").Value = 38
$ws.Range("D5").Value = 14
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 3

# Column widths: A and C:D best-fit to their widest entries, B fixed at 15
$ws.Columns.Item(1).ColumnWidth = 13.15
$ws.Columns.Item(2).ColumnWidth = 14.15
$ws.Columns.Item(3).ColumnWidth = 13.15
$ws.Columns.Item(4).ColumnWidth = 13.15

# Restore the selection left behind by the edit
$ws.Range("L14").Select()
